$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update column B (Author) -> readable "Surname et al., year" citations
# ---------------------------------------------------------------------------
$ws.Range("B2").Value  = "Verriotto et al., 2017"
$ws.Range("B3").Value  = "Eklundet al., 1996"
$ws.Range("B4").Value  = "Bajaj et al., 2011"
$ws.Range("B5").Value  = "Dianat et el.,.2013"
$ws.Range("B6").Value  = "Horne et al.,1976"
$ws.Range("B7").Value  = "Roenneberg et al.,2003"
$ws.Range("B8").Value  = "Olivier et.al.,.2016] "
$ws.Range("B9").Value  = "Buysse ei al.,1989"
$ws.Range("B10").Value = "Xie et al.,2021"
$ws.Range("B11").Value = "Wu et al.,2017 "

# ---------------------------------------------------------------------------
# 2. Update column C (Description) text that changed wording / gained counts
# ---------------------------------------------------------------------------
$ws.Range("C3").Value  = "A survey to assess electrical lighting environment in office"
$ws.Range("C8").Value  = "13 items questionnaire measuring your sleep environment quality"
$ws.Range("C10").Value = "29 Items questionnaire assessing four dimensions of biological rhythm disorder in adolescents "
$ws.Range("C11").Value = "16 dichotomous (yes/no) items questionnaire to assess ""photophobia"" and ""photophilia"""

# ---------------------------------------------------------------------------
# 3. Drop column E ("Adaptations") entirely - unused column
# ---------------------------------------------------------------------------
$ws.Columns("E").Delete()

# ---------------------------------------------------------------------------
# 4. Re-layout columns A & C: narrower, wrap text, left aligned for A10
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 26.5
$ws.Columns("C").ColumnWidth = 30.333333333333332
$ws.Columns("D").ColumnWidth = 17

$ws.Range("A1:A11").WrapText = $true
$ws.Range("C1:C11").WrapText = $true
$ws.Range("A10").HorizontalAlignment = -4131
$ws.Range("A10").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Row heights - wrapped rows grow to fit their (now narrower) columns
# ---------------------------------------------------------------------------
$ws.Rows("1").RowHeight = 17
$ws.Rows("2").RowHeight = 51
$ws.Rows("3").RowHeight = 34
$ws.Rows("4").RowHeight = 34
$ws.Rows("5").RowHeight = 34
$ws.Rows("6").RowHeight = 34
$ws.Rows("7").RowHeight = 51
$ws.Rows("8").RowHeight = 34
$ws.Rows("9").RowHeight = 34
$ws.Rows("10").RowHeight = 51
$ws.Rows("11").RowHeight = 51

# ---------------------------------------------------------------------------
# 6. Tidy up the view: scroll back to A1, select B15, maximize the window
# ---------------------------------------------------------------------------
[void]$ws.Range("A1").Select()
[void]$ws.Range("B15").Select()
$excel.ActiveWindow.WindowState = -4137
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
